# Updated cryptos list values (Price column D, Volume(1h) column E)
# Mirrors the upstream "Updated cryptos list ... with GitHub Actions" data refresh.
#
# D-column prices are locale-formatted text (dot thousands separators, e.g.
# "30.223.41", or fixed-width decimals like "12.60") so each write forces the
# literal string via a quote-prefix, then clears the resulting cell style back
# to Normal so no stray "quote prefix" number format is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.223.41"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.96%  "

$ws.Range("D3").Value = "'1.858.54"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.98%  "

$ws.Range("D4").Value = "'0.9999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").Value = "'234.31"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.46%  "

$ws.Range("D6").Value = "'0.9997"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.10%  "

$ws.Range("D7").Value = "'0.4701"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.09%  "

$ws.Range("D8").Value = "'0.2818"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.15%  "

$ws.Range("D9").Value = "'0.06552"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.55%  "

$ws.Range("E10").Value = "  +3.50%  "

$ws.Range("D11").Value = "'0.07814"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.45%  "

$ws.Range("D12").Value = "'97.04"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -7.25%  "

$ws.Range("D13").Value = "'1.864.45"
$ws.Range("D13").Style = "Normal"

$ws.Range("D14").Value = "'5.102"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.46%  "

$ws.Range("D15").Value = "'0.6649"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.65%  "

$ws.Range("D16").Value = "'283.42"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.67%  "

$ws.Range("D17").Value = "'30.260.69"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.86%  "

$ws.Range("D18").Value = "'0.9993"
$ws.Range("D18").Style = "Normal"

$ws.Range("D19").Value = "'5.438"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.27%  "

$ws.Range("D20").Value = "'12.60"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.02%  "

$ws.Range("D21").Value = "'2.106.97"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.52%  "

$ws.Range("D22").Value = "'0.000007249"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.43%  "

$ws.Range("D23").Value = "'0.9991"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.14%  "

$ws.Range("D24").Value = "'6.138"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.98%  "

$ws.Range("D25").Value = "'168.08"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.41%  "

$ws.Range("D26").Value = "'9.306"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.43%  "

$ws.Range("D27").Value = "'19.01"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.07%  "

$ws.Range("D28").Value = "'1.918"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -8.72%  "

$ws.Range("D29").Value = "'1.340"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.20%  "

$ws.Range("D30").Value = "'0.09616"
$ws.Range("D30").Style = "Normal"

$ws.Range("D31").Value = "'4.409"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.77%  "

$ws.Range("D32").Value = "'1.472"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.04%  "

$ws.Range("D33").Value = "'4.102"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.83%  "

$ws.Range("D34").Value = "'0.04675"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.64%  "

$ws.Range("D35").Value = "'1.100"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.40%  "

$ws.Range("D36").Value = "'0.7006"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.03%  "

$ws.Range("D37").Value = "'0.9987"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.13%  "

$ws.Range("D38").Value = "'2.710"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.17%  "

$ws.Range("D39").Value = "'0.01856"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.36%  "

$ws.Range("D40").Value = "'6.454"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.78%  "

$ws.Range("D41").Value = "'2.511"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.21%  "

$ws.Range("D42").Value = "'72.10"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.88%  "

$ws.Range("D43").Value = "'0.8559"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.06%  "

$ws.Range("D44").Value = "'1.937"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.70%  "

$ws.Range("D45").Value = "'104.02"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.06%  "

$ws.Range("D46").Value = "'0.4159"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.93%  "

$ws.Range("D47").Value = "'0.9991"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.06%  "

$ws.Range("D48").Value = "'1.006.35"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.19%  "

$ws.Range("D49").Value = "'7.209"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.36%  "

$ws.Range("D50").Value = "'8.996"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.61%  "

$ws.Range("D51").Value = "'33.79"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.76%  "

